$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new CURRENT_DATE column (Q) ---------------------------------
$ws.Range("Q1").Value = "CURRENT_DATE"

$ws.Range("Q2").NumberFormat = "m/d/yy;@"
$ws.Range("Q2").Formula = "=TODAY()"

$ws.Range("Q3").NumberFormat = "m/d/yy;@"
$ws.Range("Q3").Formula = "=TODAY()"

# Match column P's width (bestFit width 22) on the new column Q
$ws.Columns.Item(17).ColumnWidth = $ws.Range("P1").ColumnWidth()

# --- Extend the "Date" data validation rule to cover the new column -----
$ws.Range("P2:Q3").Validation.Delete()
$ws.Range("P2:Q3").Validation.Add(4, 1, 1, 43466, 44197)
$v = $ws.Range("P2:Q3").Validation
$v.ErrorTitle = "Date"
$v.ErrorMessage = "Valid date format e.g. 3/21/20"
$v.IgnoreBlank = $true
$v.InCellDropdown = $true
$v.ShowInput = $true
$v.ShowError = $true

# --- Update the view: scroll so column O is left-most, select Q4 --------
[void]$excel.Goto($ws.Range("O1"), $true)
[void]$ws.Range("Q4").Select()
